# New crime data collected - refresh the weekly CompStat report:
#   - bump the Volume/Number header and the reporting week date range
#   - replace the Week/28-Day/YTD/2-Year crime figures (rows 14-30)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

# --- Data table updates (rows 14-30) ---
# Row 14
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 100
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -20
$ws.Range("N14").Value = -86.666666666666

# Row 15
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 150
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = 7.692307692307
$ws.Range("I15").Value = 39
$ws.Range("J15").Value = 35
$ws.Range("K15").Value = 11.428571428571
$ws.Range("L15").Value = 8.333333333333
$ws.Range("M15").Value = 85.714285714285
$ws.Range("N15").Value = 11.428571428571

# Row 16
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 23
$ws.Range("E16").Value = 30.434782608695
$ws.Range("F16").Value = 139
$ws.Range("G16").Value = 126
$ws.Range("H16").Value = 10.31746031746
$ws.Range("I16").Value = 357
$ws.Range("J16").Value = 291
$ws.Range("K16").Value = 22.680412371134
$ws.Range("L16").Value = 86.910994764397
$ws.Range("M16").Value = -8.695652173913
$ws.Range("N16").Value = -80.888650963597

# Row 17
$ws.Range("C17").Value = 48
$ws.Range("D17").Value = 28
$ws.Range("E17").Value = 71.428571428571
$ws.Range("F17").Value = 196
$ws.Range("G17").Value = 139
$ws.Range("H17").Value = 41.007194244604
$ws.Range("I17").Value = 484
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 21
$ws.Range("L17").Value = 61.333333333333
$ws.Range("M17").Value = 91.304347826087
$ws.Range("N17").Value = 1.894736842105

# Row 18
$ws.Range("C18").Value = 34
$ws.Range("D18").Value = 31
$ws.Range("E18").Value = 9.677419354838
$ws.Range("F18").Value = 153
$ws.Range("G18").Value = 142
$ws.Range("H18").Value = 7.746478873239
$ws.Range("I18").Value = 426
$ws.Range("J18").Value = 376
$ws.Range("K18").Value = 13.297872340425
$ws.Range("L18").Value = 18.662952646239
$ws.Range("M18").Value = -25.65445026178
$ws.Range("N18").Value = -86.363636363636

# Row 19
$ws.Range("C19").Value = 148
$ws.Range("D19").Value = 117
$ws.Range("E19").Value = 26.495726495726
$ws.Range("F19").Value = 554
$ws.Range("G19").Value = 516
$ws.Range("H19").Value = 7.364341085271
$ws.Range("I19").Value = 1305
$ws.Range("J19").Value = 1514
$ws.Range("K19").Value = -13.804491413474
$ws.Range("L19").Value = 88.311688311688
$ws.Range("M19").Value = 78.522571819425
$ws.Range("N19").Value = -10.860655737704

# Row 20
$ws.Range("C20").Value = 39
$ws.Range("D20").Value = 33
$ws.Range("E20").Value = 18.181818181818
$ws.Range("F20").Value = 158
$ws.Range("G20").Value = 107
$ws.Range("H20").Value = 47.663551401869
$ws.Range("I20").Value = 420
$ws.Range("J20").Value = 277
$ws.Range("K20").Value = 51.624548736462
$ws.Range("L20").Value = 117.616580310881
$ws.Range("M20").Value = 34.185303514377
$ws.Range("N20").Value = -91.054313099041

# Row 21
$ws.Range("C21").Value = 304
$ws.Range("D21").Value = 234
$ws.Range("E21").Value = 29.914529914529
$ws.Range("F21").Value = 1216
$ws.Range("G21").Value = 1044
$ws.Range("H21").Value = 16.47509578544
$ws.Range("I21").Value = 3035
$ws.Range("J21").Value = 2899
$ws.Range("K21").Value = 4.691272852707
$ws.Range("L21").Value = 70.889639639639
$ws.Range("M21").Value = 32.706602536073
$ws.Range("N21").Value = -74.039859721153

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 5
$ws.Range("F22").Value = 31
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = 47.619047619047
$ws.Range("I22").Value = 57
$ws.Range("J22").Value = 56
$ws.Range("K22").Value = 1.785714285714
$ws.Range("L22").Value = 533.333333333333
$ws.Range("M22").Value = 39.024390243902

# Row 23
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -28.571428571428
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = -5.882352941176
$ws.Range("I23").Value = 45
$ws.Range("J23").Value = 49
$ws.Range("K23").Value = -8.163265306122
$ws.Range("L23").Value = 28.571428571428
$ws.Range("M23").Value = 50

# Row 24
$ws.Range("D24").Value = 252
$ws.Range("E24").Value = 21.825396825396
$ws.Range("F24").Value = 1190
$ws.Range("G24").Value = 1128
$ws.Range("H24").Value = 5.496453900709
$ws.Range("I24").Value = 3018
$ws.Range("J24").Value = 2761
$ws.Range("K24").Value = 9.308221658819
$ws.Range("L24").Value = 38.503900871959
$ws.Range("M24").Value = 88.154613466334

# Row 25
$ws.Range("C25").Value = 91
$ws.Range("D25").Value = 89
$ws.Range("E25").Value = 2.247191011235
$ws.Range("F25").Value = 366
$ws.Range("G25").Value = 368
$ws.Range("H25").Value = -0.543478260869
$ws.Range("I25").Value = 918
$ws.Range("J25").Value = 900
$ws.Range("K25").Value = 2
$ws.Range("L25").Value = 39.939024390243
$ws.Range("M25").Value = 3.378378378378

# Row 26
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -14.285714285714
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = -4.347826086956
$ws.Range("I26").Value = 58
$ws.Range("J26").Value = 52
$ws.Range("K26").Value = 11.538461538461
$ws.Range("L26").Value = 16

# Row 27
$ws.Range("C27").Value = 12
$ws.Range("D27").Value = 13
$ws.Range("E27").Value = -7.692307692307
$ws.Range("G27").Value = 45
$ws.Range("H27").Value = 8.888888888888
$ws.Range("I27").Value = 119
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 19
$ws.Range("L27").Value = 35.227272727272

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 12
$ws.Range("K28").Value = 50
$ws.Range("L28").Value = 9.090909090909
$ws.Range("M28").Value = 71.428571428571
$ws.Range("N28").Value = -78.181818181818

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 11
$ws.Range("K29").Value = 37.5
$ws.Range("L29").Value = 22.222222222222
$ws.Range("M29").Value = 83.333333333333
$ws.Range("N29").Value = -79.629629629629

# Row 30
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = -30
$ws.Range("I30").Value = 11
$ws.Range("J30").Value = 19
$ws.Range("K30").Value = -42.105263157894
$ws.Range("L30").Value = 10
